$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format on the Price (D) column cells we touch, so that
# numeric-looking strings (e.g. "604.15", "0.325", "1.00") are stored
# as text -- matching the original inlineStr cell type -- instead of
# being auto-coerced into numbers (which would also lose formatting
# like trailing zeros).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Cryptos list refresh: update prices and 1h volume deltas (and the
# TheGraph/Maker row swap) to the latest scrape.
$ws.Range('D2').Value = '69.684.65'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '3.500.13'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '604.15'
$ws.Range('E5').Value = '  +3.94%  '
$ws.Range('D6').Value = '171.17'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('D8').Value = '3.495.97'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  +5.63%  '
$ws.Range('D11').Value = '6.71'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('D13').Value = '46.99'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').Value = '4.070.63'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '620.27'
$ws.Range('E16').Value = '  -8.28%  '
$ws.Range('E17').Value = '  -4.03%  '
$ws.Range('D18').Value = '3.497.72'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '69.701.66'
$ws.Range('D21').Value = '17.26'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = '0.881'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').Value = '9.84'
$ws.Range('E23').Value = '  -12.07%  '
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('D25').Value = '96.16'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = '3.83'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '2.58'
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  -3.21%  '
$ws.Range('D30').Value = '33.06'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').Value = '8.37'
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('E32').Value = '  -4.83%  '
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('D34').Value = '6.93'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('D35').Value = '561.80'
$ws.Range('E35').Value = '  -5.65%  '
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('D37').Value = '3.47'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('D38').Value = '56.97'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('E39').Value = '  -4.02%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').Value = '0.141'
$ws.Range('E41').Value = '  +3.81%  '
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.323.66'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '0.325'
$ws.Range('E44').Value = '  -3.69%  '
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').Value = '32.86'
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -3.30%  '
$ws.Range('D50').Value = '134.84'
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('E51').Value = '  -1.71%  '
